# Auto-generated edit script: updates cached market-price figures
# on the "Profits" workbook, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 36684
$ws.Range("J109").Value = 36684
$ws.Range("L109").Value = 36684
$ws.Range("N109").Value = -39458
$ws.Range("H114").Value = 39372.668
$ws.Range("J114").Value = 39372.668
$ws.Range("L114").Value = 39372.668
$ws.Range("N114").Value = -48050.668
$ws.Range("H117").Value = 35519.668
$ws.Range("J117").Value = 35519.668
$ws.Range("L117").Value = 35519.668
$ws.Range("N117").Value = -44697.668
$ws.Range("H128").Value = 43037.168
$ws.Range("J128").Value = 43037.168
$ws.Range("L128").Value = 43037.168
$ws.Range("N128").Value = -52997.168

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H111").Value = 48494
$ws.Range("J111").Value = 48494
$ws.Range("L111").Value = 48494
$ws.Range("N111").Value = -56674
$ws.Range("H113").Value = 33338.5
$ws.Range("J113").Value = 33338.5
$ws.Range("L113").Value = 33338.5
$ws.Range("N113").Value = -42016.5
$ws.Range("H114").Value = 27254
$ws.Range("J114").Value = 27254
$ws.Range("L114").Value = 27254
$ws.Range("N114").Value = -35932
$ws.Range("H117").Value = 49559.75
$ws.Range("J117").Value = 49559.75
$ws.Range("L117").Value = 49559.75
$ws.Range("N117").Value = -58737.75
$ws.Range("H118").Value = 38424.25
$ws.Range("J118").Value = 38424.25
$ws.Range("L118").Value = 38424.25
$ws.Range("N118").Value = -41738.25
$ws.Range("H121").Value = 31042.4
$ws.Range("J121").Value = 31042.4
$ws.Range("L121").Value = 31042.4
$ws.Range("N121").Value = -34536.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 33782.2
$ws.Range("J112").Value = 33782.2
$ws.Range("L112").Value = 33782.2
$ws.Range("N112").Value = -36736.2
$ws.Range("H117").Value = 48935.5
$ws.Range("J117").Value = 48935.5
$ws.Range("L117").Value = 48935.5
$ws.Range("N117").Value = -58113.5
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H120").Value = 48761
$ws.Range("J120").Value = 48761
$ws.Range("L120").Value = 48761
$ws.Range("N120").Value = -58437
$ws.Range("H122").Value = 40191
$ws.Range("J122").Value = 40191
$ws.Range("L122").Value = 40191
$ws.Range("N122").Value = -49991
$ws.Range("H126").Value = 45768
$ws.Range("J126").Value = 45768
$ws.Range("L126").Value = 45768
$ws.Range("N126").Value = -55648

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 26937.4
$ws.Range("J112").Value = 26937.4
$ws.Range("L112").Value = 26937.4
$ws.Range("N112").Value = -29891.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4628.364
$ws.Range("J3").Value = 7912.6
$ws.Range("L3").Value = 23737.8
$ws.Range("N3").Value = -23961.8
$ws.Range("H24").Value = 1326
$ws.Range("I24").Value = 230
$ws.Range("J24").Value = 1600
$ws.Range("K24").Value = 690
$ws.Range("L24").Value = 4800
$ws.Range("M24").Value = -460
$ws.Range("N24").Value = -5260
$ws.Range("H113").Value = 5051.609
$ws.Range("I113").Value = 7336.8
$ws.Range("J113").Value = 766.875
$ws.Range("K113").Value = 22010.4
$ws.Range("L113").Value = 2300.625
$ws.Range("M113").Value = -19840.4
$ws.Range("N113").Value = -6640.625
$ws.Range("H131").Value = 2154.2441
$ws.Range("I131").Value = 9570.637000000001
$ws.Range("J131").Value = 1066.5067
$ws.Range("K131").Value = 28711.911
$ws.Range("L131").Value = 3199.5201
$ws.Range("M131").Value = -23671.911
$ws.Range("N131").Value = -13279.5201

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 362028.94
$ws.Range("I80").Value = 560578.3
$ws.Range("J80").Value = 4640
$ws.Range("K80").Value = 560578.3
$ws.Range("L80").Value = 4640
$ws.Range("M80").Value = -559580.3
$ws.Range("N80").Value = -6636
$ws.Range("H83").Value = 362028.94
$ws.Range("I83").Value = 560578.3
$ws.Range("J83").Value = 4640
$ws.Range("K83").Value = 2802891.5
$ws.Range("L83").Value = 23200
$ws.Range("M83").Value = -2797899.5
$ws.Range("N83").Value = -33184
$ws.Range("H99").Value = 15077.363
$ws.Range("I99").Value = 6462.75
$ws.Range("K99").Value = 6462.75
$ws.Range("M99").Value = -4216.75
$ws.Range("H110").Value = 46896
$ws.Range("J110").Value = 46896
$ws.Range("L110").Value = 46896
$ws.Range("N110").Value = -55076
$ws.Range("H114").Value = 38851.8
$ws.Range("J114").Value = 38851.8
$ws.Range("L114").Value = 38851.8
$ws.Range("N114").Value = -47529.8
$ws.Range("H116").Value = 38994.5
$ws.Range("J116").Value = 38994.5
$ws.Range("L116").Value = 38994.5
$ws.Range("N116").Value = -48172.5
$ws.Range("H119").Value = 37165
$ws.Range("J119").Value = 37165
$ws.Range("L119").Value = 37165
$ws.Range("N119").Value = -46841

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 18861
$ws.Range("J110").Value = 18861
$ws.Range("L110").Value = 18861
$ws.Range("N110").Value = -27041
$ws.Range("H116").Value = 47676
$ws.Range("J116").Value = 47676
$ws.Range("L116").Value = 47676
$ws.Range("N116").Value = -56854
$ws.Range("H117").Value = 30366.5
$ws.Range("J117").Value = 30366.5
$ws.Range("L117").Value = 30366.5
$ws.Range("N117").Value = -39544.5
$ws.Range("H120").Value = 42008
$ws.Range("J120").Value = 42008
$ws.Range("L120").Value = 42008
$ws.Range("N120").Value = -51684
$ws.Range("H121").Value = 39610.668
$ws.Range("J121").Value = 39610.668
$ws.Range("L121").Value = 39610.668
$ws.Range("N121").Value = -43104.668
